# ClosedXML DataTypes.xlsx test resource update.
#
# Commit: "Don't explicitly use Convert.ToDecimal - it leads to bounds errors."
#
# A new data row is inserted right after the existing "Double Number:" row
# (row 11) to exercise a double value that is far outside the range of a
# .NET Decimal (whose max is ~7.9E+28). Everything below the insertion
# point shifts down by one row, which Excel's row-insert does for us
# automatically (formulas/refs are not in play on this sheet, just data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 12 - pushes old rows 12..42 down to 13..43.
$ws.Rows.Item(12).Insert()

# Label in column B, matching the style of the other rows in this block.
$ws.Range("B12").Value = "Large Double Number:"

# A double value well beyond System.Decimal's bounds (~9.999E+307), stored
# as a plain number (General format, same style as its neighbours).
$largeDouble = [double]"9.999E+307"
$ws.Range("C12").Value = $largeDouble
